$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected the "Optional Delay (sec)" values in column I (rows 3-88) from
# 500 down to 60 seconds for every measured event row.
$ws.Range("I3:I88").Value = 60

# Update the sheet view: freeze the header rows (split after row 2) and set
# the active selection to E6:E7, replacing the old topLeftCell/selection.
$ws.Range("A3").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E6:E7").Select() | Out-Null
